$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell -> new value, as per the updated faturamento (billing) data
$updates = @{
    "C2"  = 50552.49
    "D2"  = 141804.95
    "E2"  = 162843.23
    "F2"  = 331318.89
    "G2"  = 406328.41

    "C3"  = 20364.19
    "D3"  = 150166.71
    "E3"  = 187284.72
    "F3"  = 343860.01
    "G3"  = 500094.77

    "C4"  = 65838.63
    "D4"  = 166620.23
    "E4"  = 193951.64
    "F4"  = 326978.27
    "G4"  = 396688.65

    "C5"  = 46495.92
    "D5"  = 130348.58
    "E5"  = 215439.97
    "F5"  = 397776.36
    "G5"  = 577277.87

    "C6"  = 47815.15
    "D6"  = 115923.13
    "E6"  = 220807.59
    "F6"  = 378480.22
    "G6"  = 108661.68

    "C7"  = 64471.91
    "D7"  = 141324.93
    "E7"  = 246540.72
    "F7"  = 404637.8

    "C8"  = 72574.36
    "D8"  = 136278.09
    "E8"  = 225495.54
    "F8"  = 390987.58

    "C9"  = 117644.99
    "D9"  = 148973.56
    "E9"  = 267994.44
    "F9"  = 382184.28

    "B10" = 40853.22
    "C10" = 114056.17
    "D10" = 153358.48
    "E10" = 292751.94
    "F10" = 388159.1

    "B11" = 33166.69
    "C11" = 120207.54
    "D11" = 134867.6
    "E11" = 332005.29
    "F11" = 418223.63

    "B12" = 18865.61
    "C12" = 129914.47
    "D12" = 151314.9
    "E12" = 249776.5
    "F12" = 310802.37

    "B13" = 23152.14
    "C13" = 106870.37
    "D13" = 201429.72
    "E13" = 251259.08
    "F13" = 384074.17
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
